$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("H 72") is removed from the data set; all rows below shift up by
# one (old row 3 -> new row 2, ..., old row 63 -> new row 62), and the
# sheet's used range shrinks from A1:F63 to A1:F62.
$ws.Rows.Item(2).Delete()
